$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 8498.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 8498.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 25494.75
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -25732.75

$ws.Range("H60").Value = 8498.25
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 8498.25
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 25494.75
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = -26462.75

$ws.Range("H132").Value = 14761.889
$ws.Range("I132").Value = 3020.5881
$ws.Range("K132").Value = 9061.764299999999
$ws.Range("M132").Value = -6531.764299999999

$ws.Range("H137").Value = 10408286
$ws.Range("I137").Value = 668703.5600000001
$ws.Range("J137").Value = 18524604
$ws.Range("K137").Value = 2006110.68
$ws.Range("L137").Value = 55573812
$ws.Range("M137").Value = -2003560.68
$ws.Range("N137").Value = -55578912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3718
$ws.Range("I45").Value = 1999.5
$ws.Range("K45").Value = 1999.5
$ws.Range("M45").Value = -1622.5

$ws.Range("H74").Value = 1213.6216
$ws.Range("I74").Value = 854.1429000000001
$ws.Range("J74").Value = 2332
$ws.Range("K74").Value = 854.1429000000001
$ws.Range("L74").Value = 2332
$ws.Range("M74").Value = 19.85709999999995
$ws.Range("N74").Value = -4080

$ws.Range("H77").Value = 1213.6216
$ws.Range("I77").Value = 854.1429000000001
$ws.Range("J77").Value = 2332
$ws.Range("K77").Value = 4270.7145
$ws.Range("L77").Value = 11660
$ws.Range("M77").Value = 97.28549999999996
$ws.Range("N77").Value = -20396

$ws.Range("H110").Value = 1136284
$ws.Range("I110").Value = 1459222.4
$ws.Range("K110").Value = 1459222.4
$ws.Range("M110").Value = -1457177.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 6945264
$ws.Range("J64").Value = 905
$ws.Range("L64").Value = 905
$ws.Range("N64").Value = -1355

$ws.Range("H67").Value = 6945264
$ws.Range("J67").Value = 905
$ws.Range("L67").Value = 905
$ws.Range("N67").Value = -2465

$ws.Range("H80").Value = 301.2
$ws.Range("J80").Value = 310.33334
$ws.Range("L80").Value = 310.33334
$ws.Range("N80").Value = -2306.33334

$ws.Range("H83").Value = 301.2
$ws.Range("J83").Value = 310.33334
$ws.Range("L83").Value = 1551.6667
$ws.Range("N83").Value = -11535.6667

$ws.Range("H107").Value = 834.5
$ws.Range("I107").Value = 821.4
$ws.Range("K107").Value = 821.4
$ws.Range("M107").Value = 1098.6

$ws.Range("H132").Value = 115635.6
$ws.Range("J132").Value = 115635.6
$ws.Range("L132").Value = 115635.6
$ws.Range("N132").Value = -125755.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 251
$ws.Range("I7").Value = 195.875
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 195.875
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -82.875
$ws.Range("N7").Value = -526

$ws.Range("H16").Value = 2763
$ws.Range("I16").Value = 2697.5
$ws.Range("K16").Value = 2697.5
$ws.Range("M16").Value = -2410.5

$ws.Range("H31").Value = 5498.9834
$ws.Range("I31").Value = 2179.5715
$ws.Range("K31").Value = 2179.5715
$ws.Range("M31").Value = -1884.5715

$ws.Range("H34").Value = 5498.9834
$ws.Range("I34").Value = 2179.5715
$ws.Range("K34").Value = 2179.5715
$ws.Range("M34").Value = -1977.5715

$ws.Range("H99").Value = 8562.071
$ws.Range("I99").Value = 3999.5
$ws.Range("J99").Value = 9322.5
$ws.Range("K99").Value = 3999.5
$ws.Range("L99").Value = 9322.5
$ws.Range("M99").Value = -2501.5
$ws.Range("N99").Value = -12318.5

$ws.Range("H105").Value = 7577090.5
$ws.Range("I105").Value = 22727272
$ws.Range("K105").Value = 22727272
$ws.Range("M105").Value = -22725525

$ws.Range("H113").Value = 2763
$ws.Range("I113").Value = 2697.5
$ws.Range("K113").Value = 2697.5
$ws.Range("M113").Value = -527.5

$ws.Range("H126").Value = 8562.071
$ws.Range("I126").Value = 3999.5
$ws.Range("J126").Value = 9322.5
$ws.Range("K126").Value = 11998.5
$ws.Range("L126").Value = 27967.5
$ws.Range("M126").Value = -9528.5
$ws.Range("N126").Value = -32907.5

$ws.Range("H132").Value = 9812904
$ws.Range("I132").Value = 11504427
$ws.Range("J132").Value = 2070.4
$ws.Range("K132").Value = 34513281
$ws.Range("L132").Value = 6211.200000000001
$ws.Range("M132").Value = -34510751
$ws.Range("N132").Value = -11271.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3784236.2
$ws.Range("J4").Value = 17507500
$ws.Range("L4").Value = 52522500
$ws.Range("N4").Value = -52522724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 108001064
$ws.Range("I3").Value = 135000080
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 135000080
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -134999964
$ws.Range("N3").Value = -5232

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = $null

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = $null

$ws.Range("H132").Value = 418296.7
$ws.Range("I132").Value = 127869.56
$ws.Range("J132").Value = 775745.4399999999
$ws.Range("K132").Value = 383608.68
$ws.Range("L132").Value = 2327236.32
$ws.Range("M132").Value = -381078.68
$ws.Range("N132").Value = -2332296.32

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5520.579
$ws.Range("I7").Value = 3089.4
$ws.Range("J7").Value = 8221.888999999999
$ws.Range("K7").Value = 3089.4
$ws.Range("L7").Value = 8221.888999999999
$ws.Range("M7").Value = -2977.4
$ws.Range("N7").Value = -8445.888999999999

$ws.Range("H40").Value = 2853.7334
$ws.Range("I40").Value = 2853.7334
$ws.Range("K40").Value = 2853.7334
$ws.Range("M40").Value = -2717.7334

$ws.Range("H126").Value = 5520.579
$ws.Range("I126").Value = 3089.4
$ws.Range("J126").Value = 8221.888999999999
$ws.Range("K126").Value = 9268.200000000001
$ws.Range("L126").Value = 24665.667
$ws.Range("M126").Value = -6798.200000000001
$ws.Range("N126").Value = -29605.667

$ws.Range("H132").Value = 3919.1936
$ws.Range("I132").Value = 3835.125
$ws.Range("J132").Value = 4207.4287
$ws.Range("K132").Value = 11505.375
$ws.Range("L132").Value = 12622.2861
$ws.Range("M132").Value = -8975.375
$ws.Range("N132").Value = -17682.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 15205.2
$ws.Range("J55").Value = 18506.75
$ws.Range("L55").Value = 18506.75
$ws.Range("N55").Value = -19060.75

$ws.Range("H107").Value = 5106.174
$ws.Range("J107").Value = 3188.625
$ws.Range("L107").Value = 9565.875
$ws.Range("N107").Value = -13405.875

$ws.Range("H126").Value = 3909.8
$ws.Range("I126").Value = 3137.25
$ws.Range("K126").Value = 9411.75
$ws.Range("M126").Value = -6941.75
